$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 58; this shifts existing rows 58..149 down to 59..150
$ws.Rows.Item(58).Insert()

# Populate the newly inserted row 58 with the new record's data
$ws.Cells.Item(58, 1).Value = 10
$ws.Cells.Item(58, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(58, 3).Value = "La Araucanía"

$ws.Cells.Item(58, 4).Value = 45272
$ws.Cells.Item(58, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(58, 5).Value = 9
$ws.Cells.Item(58, 6).Value = "Fruta"
$ws.Cells.Item(58, 7).Value = 100108
$ws.Cells.Item(58, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(58, 9).Value = 100108004
$ws.Cells.Item(58, 10).Value = "Papaya"
$ws.Cells.Item(58, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(58, 12).Value = "Primera"
$ws.Cells.Item(58, 13).Value = 80
$ws.Cells.Item(58, 14).Value = 24000
$ws.Cells.Item(58, 15).Value = 24000
$ws.Cells.Item(58, 16).Value = 24000
$ws.Cells.Item(58, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(58, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(58, 19).Value = 2400
$ws.Cells.Item(58, 20).Value = 10
